$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits -----------------------------------------------------
# Id changes
$ws.Range("A2").Value = 112139397

# Antal / Enhet / Kön (new values). "Antal" (I2) is stored as text "1"
# (not a number) in the source data, so force text the same way as the
# date cells below.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1"
$ws.Range("I2").ClearFormats()

$ws.Range("J2").Value = "ex."
$ws.Range("L2").Value = "hona"

# Dates: format cells as Text first so the date-like strings are not
# auto-converted to serial date numbers, then clear the formatting again
# so no leftover style is applied to the cell (matches plain text cells).
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2013-05-08"
$ws.Range("Y2").ClearFormats()

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2013-05-16"
$ws.Range("AA2").ClearFormats()

# --- Row 3 (new) -------------------------------------------------------
$ws.Range("A3").Value = 112156964
$ws.Range("B3").Value = 39449
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 102471
$ws.Range("F3").Value = "Åkerväddsantennmal"
$ws.Range("G3").Value = "Nemophora metallica"
$ws.Range("H3").Value = "(Poda, 1761)"
$ws.Range("M3").Value = "födosökande"
$ws.Range("P3").Value = "Åsums fure, delomr 19, 580 m NO om mc-banans ledningstorn, Sk"
$ws.Range("Q3").Value = 445828.4356342637
$ws.Range("R3").Value = 6205165.305277914
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Skåne"
$ws.Range("U3").Value = "Kristianstad"
$ws.Range("V3").Value = "Skåne"
$ws.Range("W3").Value = "Kristianstad"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2013-07-11"
$ws.Range("Y3").ClearFormats()

$ws.Range("Z3").Value = "00:00"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2013-07-11"
$ws.Range("AA3").ClearFormats()

$ws.Range("AB3").Value = "00:00"

$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AI3").Value = "i gles tallskog"
$ws.Range("AO3").Value = "på blmr av åkervädd"
$ws.Range("AW3").Value = "Nils Otto Nilsson"
$ws.Range("AX3").Value = "Nils Otto Nilsson"
$ws.Range("AY3").Value = "Åsums fure 2013"

# --- Row 4 (new) -------------------------------------------------------
$ws.Range("A4").Value = 112156959
$ws.Range("B4").Value = 39449
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 102471
$ws.Range("F4").Value = "Åkerväddsantennmal"
$ws.Range("G4").Value = "Nemophora metallica"
$ws.Range("H4").Value = "(Poda, 1761)"
$ws.Range("M4").Value = "födosökande"
$ws.Range("P4").Value = "Åsums fure, delomr 19, 580 m NO om mc-banans ledningstorn, Sk"
$ws.Range("Q4").Value = 445824.5356392039
$ws.Range("R4").Value = 6205211.776568725
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Skåne"
$ws.Range("U4").Value = "Kristianstad"
$ws.Range("V4").Value = "Skåne"
$ws.Range("W4").Value = "Kristianstad"

$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2013-07-11"
$ws.Range("Y4").ClearFormats()

$ws.Range("Z4").Value = "00:00"

$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2013-07-11"
$ws.Range("AA4").ClearFormats()

$ws.Range("AB4").Value = "00:00"

$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AI4").Value = "i tallskogsbryn"
$ws.Range("AO4").Value = "på blmr av åkervädd"
$ws.Range("AW4").Value = "Nils Otto Nilsson"
$ws.Range("AX4").Value = "Nils Otto Nilsson"
$ws.Range("AY4").Value = "Åsums fure 2013"
